$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Rename the translation keys to include the "daera_" prefix.
$ws.Range("A2").Value = "inspection_needed_export.daera_certex.heading"
$ws.Range("A3").Value = "inspection_needed_import.daera_certex.heading"

# Reset the view/selection back to the top-left corner (A1) with A3 selected.
$ws.Activate()
$ws.Range("A3").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
